$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.735.96'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").Value = '3.124.72'
$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.20'
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("E6").Value = '  +6.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '3.119.55'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +0.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.52'
$ws.Range("E10").Value = '  +1.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.153'
$ws.Range("E11").Value = '  +1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.469'
$ws.Range("E12").Value = '  -0.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000243'
$ws.Range("E13").Value = '  +0.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.71'
$ws.Range("E14").Value = '  +2.15%  '

$ws.Range("D16").Value = '3.648.59'
$ws.Range("E16").Value = '  +1.66%  '

$ws.Range("D17").Value = '67.752.17'
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.06'
$ws.Range("E18").Value = '  +0.40%  '

$ws.Range("D19").Value = '3.124.23'
$ws.Range("E19").Value = '  +1.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.52'
$ws.Range("E20").Value = '  -2.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '486.49'
$ws.Range("E21").Value = '  -0.82%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.691'
$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.74'
$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.80'
$ws.Range("E24").Value = '  +1.06%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.84'
$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  +3.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.59'
$ws.Range("E27").Value = '  +2.60%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.09'
$ws.Range("E29").Value = '  +3.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.34'
$ws.Range("E30").Value = '  +2.46%  '

$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.17'
$ws.Range("E32").Value = '  +1.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.113'
$ws.Range("E33").Value = '  +0.42%  '

$ws.Range("D34").Value = '0.0₃0953'
$ws.Range("E34").Value = '  +4.39%  '

$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '48.31'
$ws.Range("E36").Value = '  +2.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.955'
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.61'
$ws.Range("E38").Value = '  -0.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.323'
$ws.Range("E39").Value = '  +6.67%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.03'
$ws.Range("E40").Value = '  +1.99%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '49.26'
$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.124'
$ws.Range("E42").Value = '  +0.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.34'
$ws.Range("E43").Value = '  -0.24%  '

$ws.Range("E44").Value = '  +7.86%  '

$ws.Range("D45").Value = '2.795.70'
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '377.74'
$ws.Range("E46").Value = '  +0.96%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0348'
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.67'
$ws.Range("E48").Value = '  +8.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '135.64'
$ws.Range("E49").Value = '  -0.25%  '

$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.36'
$ws.Range("E51").Value = '  +8.88%  '
